$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")
$ws.Activate()

# --- Row 4 ---
$ws.Range("K4").Value = 157.654
$ws.Range("M4").Value = 157.14500000000001
$ws.Range("O4").Value = 157.14500000000001
$ws.Range("N4").Formula = "=100*(M4-G4)/M4"

# --- Row 5 ---
$ws.Range("M5").Value = 38.816000000000003
$ws.Range("O5").Value = 38.816000000000003

# --- Row 6 ---
$ws.Range("K6").Value = 37.414999999999999

# --- Row 7 ---
$ws.Range("K7").Value = 8.9920000000000009

# --- Row 8 ---
$ws.Range("K8").Value = 15.346

# --- Row 9 ---
$ws.Range("K9").Value = 3.645

# --- Row 10 ---
$ws.Range("K10").Value = 11.225
$ws.Range("M10").Value = 11.007999999999999

# Fix the shared N5:N9 formula (was using the absolute, always-row-4 $G$4
# reference; the loop-swap bug fix makes it relative to the row, like the
# other columns), then give N10 its own explicit (fixed) formula.
$ws.Range("N5:N9").Formula = "=100*(M5-G5)/M5"
$ws.Range("N10").Formula = "=100*(M10-G10)/M10"

# --- View state: move selection to O6 ---
$ws.Range("O6").Select()
